# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The workbook is a "Estado de Cuenta" (account statement) for NIT 9019155046.
# This edit:
#   1. Updates the header "Valor Mora" total (E11) and "Cant. Trabajadores" (C13).
#   2. Replaces worker CC 91537731 / ELBIS ACUÑA CACERES (period 2507) in row 16
#      with CC 1237438617 / JOSE ENRIQUE SALGADO CASSIANI (same period 2507),
#      updating the amounts for that row.
#   3. Removes the old second row for worker 91537731 (period 2506), since that
#      worker no longer has any record in the statement.
#   4. The remaining highlighted row (previously the 3rd data row, for worker
#      1237438617 / JOSE ENRIQUE SALGADO CASSIANI, period 2507) shifts up to
#      become row 17 and its period/amount are updated to the new period 2508.
#   5. All rows below shift up by one (dimension / merged cells follow).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header figures -------------------------------------------------------
# VALOR MORA total
$ws.Range("E11").Value = 129600
# Cant. Trabajadores (number of distinct workers in this statement)
$ws.Range("C13").Value = 1

# --- Data row 16: CC 91537731/ELBIS ACUÑA CACERES -> CC 1237438617/JOSE ENRIQUE SALGADO CASSIANI ---
$ws.Range("C16").Value = "1237438617"
$ws.Range("D16").Value = "JOSE ENRIQUE SALGADO CASSIANI"
$ws.Range("F16").Value = 57600
$ws.Range("G16").Value = 1800000

# --- Remove the old row 17 (CC 91537731/ELBIS ACUÑA CACERES, periodo 2506) ---
# Rows below (old row 18 onward) shift up by one.
$ws.Rows(17).Delete()

# --- New row 17 (old row 18): bump the period and the "Valor Mora" amount ---
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 72000
